$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value2 = $val
    $c.Style = "Normal"
}

Set-TextCell 'D2' '29.673.28'
Set-TextCell 'E2' '  -0.03%  '
Set-TextCell 'D3' '1.922.34'
Set-TextCell 'E3' '  -1.18%  '
Set-TextCell 'D4' '1.001'
Set-TextCell 'E4' '  -0.12%  '
Set-TextCell 'D5' '334.91'
Set-TextCell 'E5' '  -1.86%  '
Set-TextCell 'D6' '1.001'
Set-TextCell 'E6' '  -0.08%  '
Set-TextCell 'D7' '0.4657'
Set-TextCell 'E7' '  -2.57%  '
Set-TextCell 'D8' '0.4135'
Set-TextCell 'E8' '  +0.48%  '
Set-TextCell 'D9' '48.11'
Set-TextCell 'E9' '  +0.81%  '
Set-TextCell 'D10' '0.08044'
Set-TextCell 'E10' '  -2.28%  '
Set-TextCell 'D11' '1.020'
Set-TextCell 'E11' '  -1.15%  '
Set-TextCell 'D12' '22.24'
Set-TextCell 'E12' '  -1.74%  '
Set-TextCell 'D13' '1.913.82'
Set-TextCell 'E13' '  -1.25%  '
Set-TextCell 'D14' '5.996'
Set-TextCell 'E14' '  -2.25%  '
Set-TextCell 'D15' '7.174'
Set-TextCell 'E15' '  -2.38%  '
Set-TextCell 'D16' '89.62'
Set-TextCell 'E16' '  -2.21%  '
Set-TextCell 'D17' '1.001'
Set-TextCell 'E17' '  -0.21%  '
Set-TextCell 'E18' '  -1.74%  '
Set-TextCell 'D19' '0.06601'
Set-TextCell 'E19' '  -1.26%  '
Set-TextCell 'D20' '17.79'
Set-TextCell 'E20' '  -1.16%  '
Set-TextCell 'D21' '0.9973'
Set-TextCell 'E21' '  -0.30%  '
Set-TextCell 'D22' '29.663.29'
Set-TextCell 'E22' '  +0.03%  '
Set-TextCell 'E23' '  -0.63%  '
Set-TextCell 'D24' '11.60'
Set-TextCell 'E24' '  +3.43%  '
Set-TextCell 'D25' '2.205'
Set-TextCell 'E25' '  -3.65%  '
Set-TextCell 'D26' '2.126.60'
Set-TextCell 'E26' '  -2.18%  '
Set-TextCell 'D27' '157.43'
Set-TextCell 'E27' '  -2.18%  '
Set-TextCell 'D28' '19.92'
Set-TextCell 'E28' '  -1.21%  '
Set-TextCell 'D29' '2.153'
Set-TextCell 'E29' '  -0.65%  '
Set-TextCell 'D30' '5.697'
Set-TextCell 'E30' '  +1.15%  '
Set-TextCell 'D31' '117.81'
Set-TextCell 'E31' '  -4.13%  '
Set-TextCell 'D32' '1.045'
Set-TextCell 'E32' '  +3.83%  '
Set-TextCell 'D33' '0.09444'
Set-TextCell 'E33' '  -2.16%  '
Set-TextCell 'D34' '1.431'
Set-TextCell 'E34' '  -2.54%  '
Set-TextCell 'D35' '5.438'
Set-TextCell 'E35' '  -0.84%  '
Set-TextCell 'D36' '3.539'
Set-TextCell 'E36' '  -3.80%  '
Set-TextCell 'D37' '0.06152'
Set-TextCell 'E37' '  -1.60%  '
Set-TextCell 'D38' '0.02268'
Set-TextCell 'E38' '  -1.82%  '
Set-TextCell 'D39' '8.457'
Set-TextCell 'E39' '  -0.10%  '
Set-TextCell 'D40' '1.177'
Set-TextCell 'E40' '  -0.74%  '
Set-TextCell 'D41' '0.5906'
Set-TextCell 'D42' '1.001'
Set-TextCell 'E42' '  -0.12%  '
Set-TextCell 'B43' 'Algorand'
Set-TextCell 'C43' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 'D43' '0.1842'
Set-TextCell 'E43' '  -2.70%  '
Set-TextCell 'B44' 'Aptos'
Set-TextCell 'C44' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D44' '10.25'
Set-TextCell 'E44' '  -4.03%  '
Set-TextCell 'D45' '1.248'
Set-TextCell 'E45' '  -2.03%  '
Set-TextCell 'D46' '2.328'
Set-TextCell 'E46' '  -1.65%  '
Set-TextCell 'D47' '0.07537'
Set-TextCell 'E47' '  +1.60%  '
Set-TextCell 'B48' 'Decentraland'
Set-TextCell 'C48' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell 'D48' '0.5583'
Set-TextCell 'E48' '  -2.13%  '
Set-TextCell 'B49' 'EnergySwap'
Set-TextCell 'C49' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D49' '12.16'
Set-TextCell 'E49' '  -2.96%  '
Set-TextCell 'E50' '  -2.16%  '
Set-TextCell 'D51' '112.87'
Set-TextCell 'E51' '  +0.20%  '

Write-Output "Applied 104 cell updates"
